$d = $word.ActiveDocument

# Merge the split "Questions: Introduction to factorization" title runs into one.
$d.Content.Find.Execute("Questions: Introduction to factorization", $true, $false, $false, $false, $false, $true, 1, $false, "Questions: Introduction to factorization", 2)

# Merge the split "Millie Pike" author runs into one.
$d.Content.Find.Execute("Millie Pike", $true, $false, $false, $false, $false, $true, 1, $false, "Millie Pike", 2)

# Merge the split abstract sentence runs into one.
$d.Content.Find.Execute("A selection of questions for the study guide on introduction to factorization.", $true, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on introduction to factorization.", 2)
